# "Add files via upload" -- appends six new blank "Title and Content"
# slides (slide4.xml .. slide9.xml) to the end of the deck, after the
# existing 3 slides. Each new slide gets an empty Title placeholder and
# an empty Content placeholder, matching the slideLayout2.xml
# ("Title and Content") layout already used by slides 2 and 3.

$p = $ppt.ActivePresentation

$layoutTitleAndContent = 2   # ppLayoutText -> "Title and Content" (slideLayout2.xml)

# Slides 4, 5, 6, 7, 8, 9 -- appended in order at the end of the deck.
for ($n = 1; $n -le 6; $n++) {
    $newIndex = $p.Slides.Count + 1
    $slide = $p.Slides.Add($newIndex, $layoutTitleAndContent)
}
